# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates DAMSLTag (col I) and DialogAct (col J)
# values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 5;  Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 10; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 25; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 26; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 35; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 36; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 40; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 45; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 49; Tag = "qy"; Act = "Yes-No-Question" },
    @{ Row = 54; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 55; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 60; Tag = "b";  Act = "Acknowledge (Backchannel)" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.Tag
    $ws.Range("J" + $u.Row).Value = $u.Act
}
